# Adds the "SignInPage" worksheet (email look-up for the "Forgot Password"
# flow) after the existing "SignIn"/"SignUp" sheets and makes it active,
# matching the authored change.

$wb = $excel.ActiveWorkbook

$signIn = $wb.Worksheets.Item("SignIn")
$signUp = $wb.Worksheets.Item("SignUp")

# Preserve/update the cursor position saved on the two pre-existing sheets.
$signIn.Activate()
$signIn.Range("A5").Select() | Out-Null

$signUp.Activate()
$signUp.Range("C6").Select() | Out-Null

# New sheet is inserted after the last existing sheet ("SignUp").
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $signUp)
$ws.Name = "SignInPage"

# --- Column width (closest achievable match to the authored 26.53 "characters") ---
$ws.Columns.Item(1).ColumnWidth = 25.6

# --- Header row ---
$ws.Range("A1").Value = "Email Address"
$ws.Range("A1").Font.Bold = $true

# --- Email rows (with hyperlinks where the source authored them) ---
$ws.Range("A2").Value = "emailaddress@yahoo.fr"
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:emailaddress@yahoo.fr", "", "", "emailaddress@yahoo.fr") | Out-Null

$ws.Range("A3").Value = "nomane@noname.com"
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:nomane@noname.com", "", "", "nomane@noname.com") | Out-Null

$ws.Range("A4").Value = "daisy@yagm.com"
$ws.Hyperlinks.Add($ws.Range("A4"), "mailto:daisy@yagm.com", "", "", "daisy@yagm.com") | Out-Null

$ws.Range("A5").Value = "daisy@gmail.com"

$ws.Range("A6").Value = "daisy@yagm.com"
$ws.Hyperlinks.Add($ws.Range("A6"), "mailto:daisy@yagm.com", "", "", "daisy@yagm.com") | Out-Null

# Cursor position saved on the new sheet, which becomes the active tab.
$ws.Range("A6").Select() | Out-Null
$ws.Activate()

Write-Output "SignInPage sheet added"
